$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so that numeric-looking
# strings (e.g. "0.9968", "1.110") are preserved exactly as text and
# not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.305.67"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.809.56"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "0.9968"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "312.27"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "0.9966"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "0.5162"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "0.3983"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("D9").Value = "0.07896"
$ws.Range("E9").Value = "  -6.22%  "
$ws.Range("D10").Value = "1.110"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "41.13"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "6.337"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "0.9963"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "20.49"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Value = "7.326"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "1.790.25"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "92.46"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "0.00001083"
$ws.Range("E18").Value = "  -4.22%  "
$ws.Range("D19").Value = "0.06555"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "0.9961"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "17.31"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").Value = "5.988"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "28.356.43"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "11.15"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").Value = "2.220"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "160.25"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "20.58"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").Value = "2.006.97"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").Value = "2.396"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "127.51"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "0.1085"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "1.050"
$ws.Range("E32").Value = "  -4.65%  "
$ws.Range("D33").Value = "5.594"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("D34").Value = "3.650"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "0.07177"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").Value = "9.070"
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("D37").Value = "0.02333"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").Value = "0.2150"
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("D39").Value = "11.65"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "5.063"
$ws.Range("E40").Value = "  -4.31%  "
$ws.Range("D41").Value = "0.6209"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "0.9962"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("E43").Value = "  -3.35%  "
$ws.Range("D44").Value = "13.26"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("D45").Value = "1.325"
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("D46").Value = "0.5988"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").Value = "3.743"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").Value = "126.04"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").Value = "1.214"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "1.937"
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("D51").Value = "0.06865"
$ws.Range("E51").Value = "  -1.90%  "

# Remove the temporary text-number-format styling so the cells keep
# their original (default) style, now that the text values are locked in.
$ws.Range("D2:E51").ClearFormats()
